$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.03509803921568627, 0.4809803921568628, 0.4637941176470588, 0.6810243150189711, 0.6868868103232719, 51)
    3  = @(-0.1432, 0.5379999999999999, 0.5124119999999999, 0.7158295886592003, 0.7084805241227344, 50)
    4  = @(-0.03999999999999999, 0.4955102040816328, 0.456461224489796, 0.6756191415951712, 0.6814231431350127, 49)
    5  = @(-0.1310416666666667, 0.5352083333333333, 0.4989104166666666, 0.7063359092292184, 0.7014187433666847, 48)
    6  = @(-0.03829787234042552, 0.4770212765957446, 0.4602893617021277, 0.6784462850529346, 0.6846875490163257, 47)
    7  = @(-0.1304347826086956, 0.5034782608695652, 0.4911652173913044, 0.7008318039239546, 0.6961958744866088, 46)
    8  = @(-0.04133333333333333, 0.4848888888888888, 0.4783555555555555, 0.691632529278052, 0.698197679744068, 45)
    9  = @(-0.1425, 0.5220454545454546, 0.510906818181818, 0.7147774606000235, 0.7085265582325974, 44)
    10 = @(-0.1113953488372093, 0.4225581395348837, 0.3647790697674418, 0.6039694278417094, 0.6006329256762595, 43)
    11 = @(-0.1119047619047619, 0.4828571428571428, 0.4498904761904762, 0.6707387540544203, 0.669354396251828, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
